# "SQL queries for category improved for percent speed"
# On slide 2 ("About us"), the "Hobbies" paragraph's single run is split
# into three runs (no visible formatting change), and the accidental
# double space before "camping" is corrected to a single space:
#   "Hobbies: guitars, karaoke,  camping etc."
# becomes three runs:
#   "Hobbies: guitars, " + "karaoke, " + "camping etc."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Find the shape/paragraph that actually holds the "Hobbies: ..." line
# instead of assuming fixed shape/paragraph indices.
$tr = $null
$hobbiesIndex = -1
for ($si = 1; $si -le $s.Shapes.Count; $si++) {
    $shp = $s.Shapes.Item($si)
    if (-not $shp.HasTextFrame) { continue }
    $candidateTr = $shp.TextFrame.TextRange
    for ($i = 1; $i -le $candidateTr.Paragraphs().Count; $i++) {
        if ($candidateTr.Paragraphs($i).Text.StartsWith("Hobbies:")) {
            $tr = $candidateTr
            $hobbiesIndex = $i
            break
        }
    }
    if ($hobbiesIndex -ne -1) { break }
}

$para = $tr.Paragraphs($hobbiesIndex)
$run1 = $para.Runs(1)
$start = $run1.Start
$fullText = $run1.Text

$karaokeAt = $fullText.IndexOf("karaoke")
$campingAt = $fullText.IndexOf("camping")
$midLen = $campingAt - $karaokeAt

# Replacing this middle slice (which also absorbs the stray double space)
# with its corrected, single-spaced text forces the engine to split the
# single run into three sibling runs without touching visible formatting.
$midRange = $tr.Characters($start + $karaokeAt, $midLen)
$midRange.Text = "karaoke, "
